$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5416.5
$ws.Range("I40").Value = 4166.6665
$ws.Range("J40").Value = 5833.1113
$ws.Range("K40").Value = 4166.6665
$ws.Range("L40").Value = 5833.1113
$ws.Range("M40").Value = -3991.6665
$ws.Range("N40").Value = -6183.1113
$ws.Range("H62").Value = 35178
$ws.Range("I62").Value = 40282.668
$ws.Range("J62").Value = 33263.75
$ws.Range("K62").Value = 40282.668
$ws.Range("L62").Value = 33263.75
$ws.Range("M62").Value = -39658.668
$ws.Range("N62").Value = -34511.75
$ws.Range("H65").Value = 35178
$ws.Range("I65").Value = 40282.668
$ws.Range("J65").Value = 33263.75
$ws.Range("K65").Value = 201413.34
$ws.Range("L65").Value = 166318.75
$ws.Range("M65").Value = -198293.34
$ws.Range("N65").Value = -172558.75
$ws.Range("H116").Value = 37871.08
$ws.Range("I116").Value = 15200.8
$ws.Range("K116").Value = 15200.8
$ws.Range("M116").Value = -11758.8
$ws.Range("H125").Value = 36285.5
$ws.Range("I125").Value = 47625.223
$ws.Range("K125").Value = 428627.007
$ws.Range("M125").Value = -426167.007
$ws.Range("H137").Value = 4030.3062
$ws.Range("I137").Value = 1887.7742
$ws.Range("K137").Value = 5663.3226
$ws.Range("M137").Value = -3113.3226
$ws.Range("H138").Value = 2256.2778
$ws.Range("I138").Value = 1355.2
$ws.Range("J138").Value = 3382.625
$ws.Range("K138").Value = 4065.6
$ws.Range("L138").Value = 10147.875
$ws.Range("M138").Value = 1074.4
$ws.Range("N138").Value = -20427.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2412.9285
$ws.Range("I2").Value = 1898.8334
$ws.Range("K2").Value = 1898.8334
$ws.Range("M2").Value = -1785.8334
$ws.Range("H32").Value = 3645114.5
$ws.Range("I32").Value = 719070.4399999999
$ws.Range("K32").Value = 719070.4399999999
$ws.Range("M32").Value = -718783.4399999999
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").Value = ""
$ws.Range("H74").Value = 2321.8076
$ws.Range("I74").Value = 1376.2
$ws.Range("K74").Value = 1376.2
$ws.Range("M74").Value = -502.2
$ws.Range("H77").Value = 2321.8076
$ws.Range("I77").Value = 1376.2
$ws.Range("K77").Value = 6881
$ws.Range("M77").Value = -2513
$ws.Range("H116").Value = 2412.9285
$ws.Range("I116").Value = 1898.8334
$ws.Range("K116").Value = 1898.8334
$ws.Range("M116").Value = 395.1666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2412.9285
$ws.Range("I3").Value = 1898.8334
$ws.Range("K3").Value = 1898.8334
$ws.Range("M3").Value = -1784.8334
$ws.Range("H105").Value = 18231.172
$ws.Range("I105").Value = 4400.6665
$ws.Range("J105").Value = 54536.25
$ws.Range("K105").Value = 4400.6665
$ws.Range("L105").Value = 54536.25
$ws.Range("M105").Value = -2653.6665
$ws.Range("N105").Value = -58030.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8122.25
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").Value = ""
$ws.Range("H34").Value = 8122.25
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").Value = ""
$ws.Range("H92").Value = 27067.334
$ws.Range("J92").Value = 27067.334
$ws.Range("L92").Value = 27067.334
$ws.Range("N92").Value = -32059.334
$ws.Range("H99").Value = 20255.37
$ws.Range("J99").Value = 11653.223
$ws.Range("L99").Value = 11653.223
$ws.Range("N99").Value = -14649.223
$ws.Range("H126").Value = 20255.37
$ws.Range("J126").Value = 11653.223
$ws.Range("L126").Value = 34959.669
$ws.Range("N126").Value = -39899.669

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 425.66666
$ws.Range("I12").Value = 889
$ws.Range("J12").Value = 383.54544
$ws.Range("K12").Value = 2667
$ws.Range("L12").Value = 1150.63632
$ws.Range("M12").Value = -2494
$ws.Range("N12").Value = -1496.63632
$ws.Range("H35").Value = 201.875
$ws.Range("I35").Value = 95.833336
$ws.Range("J35").Value = 520
$ws.Range("K35").Value = 287.500008
$ws.Range("L35").Value = 1560
$ws.Range("M35").Value = 0.4999920000000202
$ws.Range("N35").Value = -2136
$ws.Range("H37").Value = 125086870
$ws.Range("J37").Value = 125086870
$ws.Range("L37").Value = 375260610
$ws.Range("N37").Value = -375260834
$ws.Range("H38").Value = 83.15000000000001
$ws.Range("I38").Value = 21.083334
$ws.Range("K38").Value = 63.250002
$ws.Range("M38").Value = 283.749998
$ws.Range("H50").Value = 59272.766
$ws.Range("I50").Value = 125134.75
$ws.Range("J50").Value = 728.7778
$ws.Range("K50").Value = 375404.25
$ws.Range("L50").Value = 2186.3334
$ws.Range("M50").Value = -374923.25
$ws.Range("N50").Value = -3148.3334
$ws.Range("H53").Value = 59272.766
$ws.Range("I53").Value = 125134.75
$ws.Range("J53").Value = 728.7778
$ws.Range("K53").Value = 375404.25
$ws.Range("L53").Value = 2186.3334
$ws.Range("M53").Value = -374923.25
$ws.Range("N53").Value = -3148.3334
$ws.Range("H68").Value = 1282.2858
$ws.Range("J68").Value = 1244.3572
$ws.Range("L68").Value = 3733.0716
$ws.Range("N68").Value = -5355.071599999999
$ws.Range("H71").Value = 1282.2858
$ws.Range("J71").Value = 1244.3572
$ws.Range("L71").Value = 11199.2148
$ws.Range("N71").Value = -19311.2148
$ws.Range("H94").Value = 152692.58
$ws.Range("I94").Value = 500474.5
$ws.Range("J94").Value = 13579.8
$ws.Range("K94").Value = 1501423.5
$ws.Range("L94").Value = 40739.39999999999
$ws.Range("M94").Value = -1500747.5
$ws.Range("N94").Value = -42091.39999999999
$ws.Range("H97").Value = 444.7143
$ws.Range("J97").Value = 374.66666
$ws.Range("L97").Value = 1123.99998
$ws.Range("N97").Value = -2115.99998
$ws.Range("H103").Value = 27781568
$ws.Range("J103").Value = 30306254
$ws.Range("L103").Value = 90918762
$ws.Range("N103").Value = -90920520
$ws.Range("H113").Value = 899.5
$ws.Range("I113").Value = 947.5
$ws.Range("J113").Value = 851.5
$ws.Range("K113").Value = 2842.5
$ws.Range("L113").Value = 2554.5
$ws.Range("M113").Value = -672.5
$ws.Range("N113").Value = -6894.5
$ws.Range("H114").Value = 11112080
$ws.Range("I114").Value = 22223136
$ws.Range("J114").Value = 1025.4445
$ws.Range("K114").Value = 66669408
$ws.Range("L114").Value = 3076.3335
$ws.Range("M114").Value = -66666154
$ws.Range("N114").Value = -9584.333500000001
$ws.Range("H128").Value = 385287.72
$ws.Range("I128").Value = 385287.72
$ws.Range("K128").Value = 1155863.16
$ws.Range("M128").Value = -1150883.16
$ws.Range("H131").Value = 25117640
$ws.Range("I131").Value = 38539396
$ws.Range("J131").Value = 191517.28
$ws.Range("K131").Value = 115618188
$ws.Range("L131").Value = 574551.84
$ws.Range("M131").Value = -115613148
$ws.Range("N131").Value = -584631.84
$ws.Range("H137").Value = 8336954.5
$ws.Range("J137").Value = 7483.25
$ws.Range("L137").Value = 22449.75
$ws.Range("N137").Value = -32649.75
$ws.Range("H140").Value = 3649.3333
$ws.Range("I140").Value = 3131.7273
$ws.Range("J140").Value = 4462.7144
$ws.Range("K140").Value = 9395.1819
$ws.Range("L140").Value = 13388.1432
$ws.Range("M140").Value = -4215.1819
$ws.Range("N140").Value = -23748.1432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 1298.5
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").Value = ""
$ws.Range("H37").Value = 1298.5
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").Value = ""
$ws.Range("H126").Value = 51304.08
$ws.Range("I126").Value = 2635.3333
$ws.Range("K126").Value = 7905.999899999999
$ws.Range("M126").Value = -5435.999899999999
$ws.Range("H132").Value = 7526.4287
$ws.Range("I132").Value = 9138.1
$ws.Range("J132").Value = 3497.25
$ws.Range("K132").Value = 27414.3
$ws.Range("L132").Value = 10491.75
$ws.Range("M132").Value = -24884.3
$ws.Range("N132").Value = -15551.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 15030.917
$ws.Range("I7").Value = 12863.444
$ws.Range("K7").Value = 12863.444
$ws.Range("M7").Value = -12751.444
$ws.Range("H16").Value = 1167.8667
$ws.Range("I16").Value = 1101.3846
$ws.Range("K16").Value = 1101.3846
$ws.Range("M16").Value = -931.3846000000001
$ws.Range("H22").Value = 1325.3914
$ws.Range("I22").Value = 948.25
$ws.Range("J22").Value = 1736.8182
$ws.Range("K22").Value = 948.25
$ws.Range("L22").Value = 1736.8182
$ws.Range("M22").Value = -653.25
$ws.Range("N22").Value = -2326.8182
$ws.Range("H27").Value = 1325.3914
$ws.Range("I27").Value = 948.25
$ws.Range("J27").Value = 1736.8182
$ws.Range("K27").Value = 948.25
$ws.Range("L27").Value = 1736.8182
$ws.Range("M27").Value = -841.25
$ws.Range("N27").Value = -1950.8182
$ws.Range("H40").Value = 5320.533
$ws.Range("I40").Value = 5129.5
$ws.Range("K40").Value = 5129.5
$ws.Range("M40").Value = -4993.5
$ws.Range("H46").Value = 33759
$ws.Range("I46").Value = 43426.7
$ws.Range("J46").Value = 1533.3334
$ws.Range("K46").Value = 43426.7
$ws.Range("L46").Value = 1533.3334
$ws.Range("M46").Value = -43238.7
$ws.Range("N46").Value = -1909.3334
$ws.Range("H126").Value = 15030.917
$ws.Range("I126").Value = 12863.444
$ws.Range("K126").Value = 38590.33199999999
$ws.Range("M126").Value = -36120.33199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2659.3333
$ws.Range("I126").Value = 2315.4375
$ws.Range("K126").Value = 6946.3125
$ws.Range("M126").Value = -4476.3125
